$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignment 05")

# Self-evaluation score entered for "Time devoted to the Topical Assignment" (hours)
$ws.Range("D2").Value = 3

# Topical Assignment checklist: first section (rows 8-30) answered "Yes" instead of "No"
$ws.Range("D8:D30").Value = "Yes"

# Reflect the reviewer's current scroll/zoom position on the sheet
$ws.Application.ActiveWindow.Zoom = 90
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D33").Select()
